$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ENGLISH")
$ws2 = $wb.Worksheets.Item("NOTES")

# --- ENGLISH sheet: insert new row 68 ("extricate") and shift the rest down ---
$ws1.Rows.Item(68).Insert()
$ws1.Range("A68").Value = 'extricate'
$ws1.Range("B68").Value = 'free from a constraint or difficulty'
$ws1.Range("F68").Value = '2021-11-20 14:57:27.485397'
$ws1.Range("E68").Value = 0

# --- ENGLISH sheet: append 9 new vocabulary rows at the end (rows 136-144) ---
$ws1.Range("A136").Value = 'restrained'
$ws1.Range("C136").Value = 'self-controlled'
$ws1.Range("F136").Value = '2021-11-20 15:03:49.858279'
$ws1.Range("E136").Value = 0

$ws1.Range("A137").Value = 'condense'
$ws1.Range("B137").Value = 'make denser or more concentrated'
$ws1.Range("F137").Value = '2021-11-20 15:04:42.212772'
$ws1.Range("E137").Value = 0

$ws1.Range("A138").Value = 'apt'
$ws1.Range("C138").Value = 'inclined;suitable'
$ws1.Range("F138").Value = '2021-11-20 15:05:21.928371'
$ws1.Range("E138").Value = 0

$ws1.Range("A139").Value = 'sentiment'
$ws1.Range("B139").Value = 'a view or opinion that is held or expressed'
$ws1.Range("C139").Value = 'view;feeling'
$ws1.Range("F139").Value = '2021-11-20 15:06:18.444516'
$ws1.Range("E139").Value = 0

$ws1.Range("A140").Value = 'fuming'
$ws1.Range("B140").Value = 'feeling, showing or expressing great anger'
$ws1.Range("F140").Value = '2021-11-20 15:07:39.355342'
$ws1.Range("E140").Value = 0

$ws1.Range("A141").Value = 'insidiously'
$ws1.Range("B141").Value = 'in a gradual, subtle way, but with harmful effects'
$ws1.Range("F141").Value = '2021-11-20 15:08:28.761904'
$ws1.Range("E141").Value = 0

$ws1.Range("A142").Value = 'fret'
$ws1.Range("B142").Value = 'be constantly or visibly anxious'
$ws1.Range("C142").Value = 'worry;trouble'
$ws1.Range("F142").Value = '2021-11-20 15:10:17.940034'
$ws1.Range("E142").Value = 0

$ws1.Range("A143").Value = 'unsolicited'
$ws1.Range("B143").Value = 'not asked for; given or done voluntarily'
$ws1.Range("C143").Value = 'uninvited'
$ws1.Range("F143").Value = '2021-11-20 15:12:00.23661'
$ws1.Range("E143").Value = 0

$ws1.Range("A144").Value = 'resemble'
$ws1.Range("C144").Value = 'look like'
$ws1.Range("F144").Value = '2021-11-20 15:12:28.779915'
$ws1.Range("E144").Value = 0

# --- NOTES sheet: append 7 new rows at the end (rows 29-35) ---
$ws2.Range("A29").Value = 'Get over the fear of waste'
$ws2.Range("B29").Value = 'essentialism'

$ws2.Range("A30").Value = 'Stop making casual commitments'
$ws2.Range("B30").Value = 'essentialism'

$ws2.Range("A31").Value = 'Pause before you speak'

$ws2.Range("A32").Value = 'get over the fear of missing out'

$ws2.Range("A33").Value = 'I saw the angel in the marble and carved until I set him free'

$ws2.Range("A34").Value = 'No is a complete sentence'

$ws2.Range("A35").Value = 'If you don''t set boundaries - there won''t be any'

